$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.530.79'
$ws.Range('E2').Value = '  -4.42%  '
$ws.Range('D3').Value = '1.807.48'
$ws.Range('E3').Value = '  -3.18%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').Value = "'274.80"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -8.52%  '
$ws.Range('D6').Value = "'1.003"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.18%  '
$ws.Range('D7').Value = "'0.4987"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -6.18%  '
$ws.Range('D8').Value = "'0.3424"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -8.04%  '
$ws.Range('D9').Value = "'44.08"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -2.73%  '
$ws.Range('D10').Value = "'0.06624"
$ws.Range('D10').Style = 'Normal'
$ws.Range('D11').Value = "'19.39"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -9.55%  '
$ws.Range('D12').Value = "'0.7949"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -10.30%  '
$ws.Range('D13').Value = "'0.07843"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -4.14%  '
$ws.Range('D14').Value = '1.806.33'
$ws.Range('E14').Value = '  -3.38%  '
$ws.Range('D15').Value = "'5.010"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -5.17%  '
$ws.Range('D16').Value = "'86.22"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -6.47%  '
$ws.Range('E17').Value = '  +0.02%  '
$ws.Range('D18').Value = "'13.91"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -5.91%  '
$ws.Range('D19').Value = "'1.002"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.13%  '
$ws.Range('D20').Value = "'0.000007915"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -6.61%  '
$ws.Range('D21').Value = '25.571.74'
$ws.Range('E21').Value = '  -4.46%  '
$ws.Range('D22').Value = "'4.697"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -5.42%  '
$ws.Range('D23').Value = "'9.853"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -6.96%  '
$ws.Range('D24').Value = "'6.077"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -4.36%  '
$ws.Range('D25').Value = "'2.235"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -2.30%  '
$ws.Range('D26').Value = "'142.66"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.26%  '
$ws.Range('D27').Value = "'1.657"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -3.65%  '
$ws.Range('D28').Value = "'17.00"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -5.50%  '
$ws.Range('D29').Value = "'108.12"
$ws.Range('D29').Style = 'Normal'
$ws.Range('D30').Value = "'4.235"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -9.51%  '
$ws.Range('D31').Value = "'4.175"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -9.52%  '
$ws.Range('D32').Value = "'0.08700"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -4.55%  '
$ws.Range('D33').Value = "'0.04767"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -4.71%  '
$ws.Range('B34').Value = 'ARBITRUM'
$ws.Range('C34').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D34').Value = "'1.121"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -4.11%  '
$ws.Range('B35').Value = 'ImmutableX'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D35').Value = "'0.7112"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -11.05%  '
$ws.Range('B36').Value = 'HuobiToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D36').Value = "'2.845"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -3.36%  '
$ws.Range('E37').Value = '  +0.07%  '
$ws.Range('D38').Value = "'3.115"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.89%  '
$ws.Range('D39').Value = "'2.303"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -14.06%  '
$ws.Range('D40').Value = "'0.01822"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -6.10%  '
$ws.Range('D41').Value = "'0.4977"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -17.85%  '
$ws.Range('D42').Value = "'0.9332"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -12.07%  '
$ws.Range('D43').Value = "'114.71"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.08%  '
$ws.Range('D44').Value = "'6.133"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -5.31%  '
$ws.Range('D45').Value = "'1.002"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.21%  '
$ws.Range('D46').Value = "'7.701"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -11.55%  '
$ws.Range('D47').Value = "'0.1341"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -9.93%  '
$ws.Range('D48').Value = "'0.4346"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -16.02%  '
$ws.Range('D49').Value = "'35.96"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -3.70%  '
$ws.Range('D50').Value = "'9.112"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -8.04%  '
$ws.Range('D51').Value = "'0.05806"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -4.25%  '
